$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.953.40'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '1.638.11'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").Value = "'215.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("D9").Value = "'0.0638"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = "'19.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.20%  '
$ws.Range("D11").Value = "'0.0796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("D12").Value = '1.864.97'
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").Value = '1.633.56'
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").Value = '0.0₃0764'
$ws.Range("E16").Value = '  -0.44%  '
$ws.Range("D17").Value = "'63.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.91%  '
$ws.Range("D18").Value = '26.021.67'
$ws.Range("E18").Value = '  -1.05%  '
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = "'193.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.85%  '
$ws.Range("E21").Value = '  -1.92%  '
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").Value = "'144.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("E27").Value = '  +3.34%  '
$ws.Range("D28").Value = "'6.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("E34").Value = '  -4.70%  '
$ws.Range("E35").Value = '  +1.76%  '
$ws.Range("D36").Value = "'0.901"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.68%  '
$ws.Range("D37").Value = '1.137.17'
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("E42").Value = '  -3.13%  '
$ws.Range("D43").Value = "'99.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("D44").Value = "'0.798"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").Value = '1.774.69'
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("E46").Value = '  +4.24%  '
$ws.Range("D47").Value = "'56.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").Value = "'0.0531"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.90%  '
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").Value = "'7.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("E51").Value = '  -0.59%  '
